$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting of the existing header row
# (G1) so it matches the other header cells, then set its value/text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add data values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
